# Update "want to go" counters (column F) across the four worksheets to
# reflect the latest scraped snapshot, matching commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 373
$ws.Range("F5").Value = 1280
$ws.Range("F6").Value = 214
$ws.Range("F7").Value = 2452
$ws.Range("F8").Value = 870
$ws.Range("F9").Value = 18468
$ws.Range("F11").Value = 1862
$ws.Range("F12").Value = 651
$ws.Range("F13").Value = 594
$ws.Range("F14").Value = 317
$ws.Range("F15").Value = 592
$ws.Range("F17").Value = 193
$ws.Range("F18").Value = 65
$ws.Range("F19").Value = 315
$ws.Range("F21").Value = 92

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 18
$ws.Range("F6").Value = 27
$ws.Range("F8").Value = 116
$ws.Range("F13").Value = 8
$ws.Range("F21").Value = 20

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5855
$ws.Range("F3").Value = 553

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5855
$ws.Range("F4").Value = 553
$ws.Range("F6").Value = 373
$ws.Range("F7").Value = 18
$ws.Range("F10").Value = 1280
$ws.Range("F12").Value = 214
$ws.Range("F14").Value = 27
$ws.Range("F15").Value = 2452
$ws.Range("F16").Value = 870
$ws.Range("F17").Value = 18468
$ws.Range("F20").Value = 116
$ws.Range("F21").Value = 116
$ws.Range("F22").Value = 1862
$ws.Range("F23").Value = 651
$ws.Range("F25").Value = 594
$ws.Range("F26").Value = 317
$ws.Range("F27").Value = 592
$ws.Range("F29").Value = 193
$ws.Range("F31").Value = 65
$ws.Range("F34").Value = 315
$ws.Range("F35").Value = 8
$ws.Range("F39").Value = 92
$ws.Range("F46").Value = 20
